# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch count),
# writing the newly calculated s_vals into column G (K) for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 3
    11 = 1
    12 = 5
    13 = 2
    14 = 4
    15 = 5
    16 = 2
    17 = 1
    19 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
